# Natmi following Dr Hou advice
# Updates the Cd14-Itgb1 ligand-receptor pair sheet: ligand/receptor-expressing
# cell counts move from 1 to 3 replicates, and all dependent expression /
# specificity statistics (avg/total expression, specificity scores, edge
# weights) are refreshed to the recomputed values for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.029662
$ws.Range("H2").Value = 3.088986
$ws.Range("I2").Value = 0.008285277389098726
$ws.Range("J2").Value = 0.008285277389098728
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 158.0829950454747
$ws.Range("R2").Value = 1422.746955409272
$ws.Range("S2").Value = 0.002628261467229324
$ws.Range("T2").Value = 0.002628261467229324
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.029662
$ws.Range("H3").Value = 3.088986
$ws.Range("I3").Value = 0.008285277389098726
$ws.Range("J3").Value = 0.008285277389098728
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 173.8066394471654
$ws.Range("R3").Value = 1564.259755024488
$ws.Range("S3").Value = 0.002889680152354133
$ws.Range("T3").Value = 0.002889680152354134
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.029662
$ws.Range("H4").Value = 3.088986
$ws.Range("I4").Value = 0.008285277389098726
$ws.Range("J4").Value = 0.008285277389098728
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 70.11001850404666
$ws.Range("R4").Value = 630.99016653642
$ws.Range("S4").Value = 0.001165637455489211
$ws.Range("T4").Value = 0.001165637455489212
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.029662
$ws.Range("H5").Value = 3.088986
$ws.Range("I5").Value = 0.008285277389098726
$ws.Range("J5").Value = 0.008285277389098728
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 96.33792900652601
$ws.Range("R5").Value = 867.0413610587341
$ws.Range("S5").Value = 0.001601698314026058
$ws.Range("T5").Value = 0.001601698314026058
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.017189
$ws.Range("H6").Value = 3.051567
$ws.Range("I6").Value = 0.00818491215771772
$ws.Range("J6").Value = 0.008184912157717722
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 156.168027612276
$ws.Range("R6").Value = 1405.512248510484
$ws.Range("S6").Value = 0.002596423538587933
$ws.Range("T6").Value = 0.002596423538587933
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.017189
$ws.Range("H7").Value = 3.051567
$ws.Range("I7").Value = 0.00818491215771772
$ws.Range("J7").Value = 0.008184912157717722
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 171.701200755804
$ws.Range("R7").Value = 1545.310806802236
$ws.Range("S7").Value = 0.002854675480393515
$ws.Range("T7").Value = 0.002854675480393516
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.017189
$ws.Range("H8").Value = 3.051567
$ws.Range("I8").Value = 0.00818491215771772
$ws.Range("J8").Value = 0.008184912157717722
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 69.26072790110999
$ws.Range("R8").Value = 623.3465511099899
$ws.Range("S8").Value = 0.001151517291802179
$ws.Range("T8").Value = 0.00115151729180218
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.017189
$ws.Range("H9").Value = 3.051567
$ws.Range("I9").Value = 0.00818491215771772
$ws.Range("J9").Value = 0.008184912157717722
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 95.170921786197
$ws.Range("R9").Value = 856.538296075773
$ws.Range("S9").Value = 0.001582295846934093
$ws.Range("T9").Value = 0.001582295846934093
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 120.9972243333333
$ws.Range("H10").Value = 362.991673
$ws.Range("I10").Value = 0.9736161642487271
$ws.Range("J10").Value = 0.9736161642487271
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 18576.58495195755
$ws.Range("R10").Value = 167189.264567618
$ws.Range("S10").Value = 0.3088511981184139
$ws.Range("T10").Value = 0.3088511981184139
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 120.9972243333333
$ws.Range("H11").Value = 362.991673
$ws.Range("I11").Value = 0.9736161642487271
$ws.Range("J11").Value = 0.9736161642487271
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 20424.29549095863
$ws.Range("R11").Value = 183818.6594186277
$ws.Range("S11").Value = 0.3395709248724085
$ws.Range("T11").Value = 0.3395709248724085
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 120.9972243333333
$ws.Range("H12").Value = 362.991673
$ws.Range("I12").Value = 0.9736161642487271
$ws.Range("J12").Value = 0.9736161642487271
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 8238.740127292534
$ws.Range("R12").Value = 74148.66114563281
$ws.Range("S12").Value = 0.1369759170418681
$ws.Range("T12").Value = 0.1369759170418681
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 120.9972243333333
$ws.Range("H13").Value = 362.991673
$ws.Range("I13").Value = 0.9736161642487271
$ws.Range("J13").Value = 0.9736161642487271
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 11320.82373420731
$ws.Range("R13").Value = 101887.4136078658
$ws.Range("S13").Value = 0.1882181242160366
$ws.Range("T13").Value = 0.1882181242160366
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.232029333333333
$ws.Range("H14").Value = 3.696088
$ws.Range("I14").Value = 0.009913646204456457
$ws.Range("J14").Value = 0.009913646204456457
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 189.1522528725084
$ws.Range("R14").Value = 1702.370275852576
$ws.Range("S14").Value = 0.003144813757617773
$ws.Range("T14").Value = 0.003144813757617774
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.232029333333333
$ws.Range("H15").Value = 3.696088
$ws.Range("I15").Value = 0.009913646204456457
$ws.Range("J15").Value = 0.009913646204456457
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 207.9661851432783
$ws.Range("R15").Value = 1871.695666289504
$ws.Range("S15").Value = 0.003457611052608941
$ws.Range("T15").Value = 0.003457611052608941
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.232029333333333
$ws.Range("H16").Value = 3.696088
$ws.Range("I16").Value = 0.009913646204456457
$ws.Range("J16").Value = 0.009913646204456457
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 83.88927566281778
$ws.Range("R16").Value = 755.00348096536
$ws.Range("S16").Value = 0.001394729083130907
$ws.Range("T16").Value = 0.001394729083130908
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.232029333333333
$ws.Range("H17").Value = 3.696088
$ws.Range("I17").Value = 0.009913646204456457
$ws.Range("J17").Value = 0.009913646204456457
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 115.2719576410747
$ws.Range("R17").Value = 1037.447618769672
$ws.Range("S17").Value = 0.001916492311098835
$ws.Range("T17").Value = 0.001916492311098835
